$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 800
$ws.Range("D3").Value = 800
$ws.Range("D4").Value = 600
$ws.Range("D6").Value = 800
$ws.Range("D7").Value = 800
$ws.Range("D12").Value = 600
$ws.Range("D13").Value = 550
$ws.Range("D15").Value = 664.29
$ws.Range("D16").Value = 800
$ws.Range("D17").Value = 550
$ws.Range("D19").Value = 600
$ws.Range("D21").Value = 600
$ws.Range("D22").Value = 800
$ws.Range("D23").Value = 600
$ws.Range("D24").Value = 600
